# Login Tests with Select
# - Rename "SecondTest" -> "addCustomerTest"
# - Add postcode column (C) + a 4th data row (Tim Fisher) to addCustomerTest
# - Add a new sheet "openAccountTest" with customer/currency data
# - Move the active/selected tab to the new last sheet

$wb = $excel.ActiveWorkbook

# --- Rename the second sheet ---
$wsCustomer = $wb.Worksheets.Item(2)
$wsCustomer.Name = "addCustomerTest"

# --- Add the postcode column header + new row of data ---
$wsCustomer.Range("C1").Value = "postcode"
$wsCustomer.Range("C2").Value = 12345
$wsCustomer.Range("C3").Value = 56789

$wsCustomer.Range("A4").Value = "Tim"
$wsCustomer.Range("B4").Value = "Fisher"
$wsCustomer.Range("C4").Value = 45678

# addCustomerTest is no longer the selected/active tab; update its selection
[void]$wsCustomer.Range("G8").Select()

# --- Add the new "openAccountTest" sheet after addCustomerTest ---
$wsAccount = $wb.Worksheets.Add($null, $wsCustomer)
$wsAccount.Name = "openAccountTest"

$wsAccount.Range("A1").Value = "customer"
$wsAccount.Range("B1").Value = "currency"
$wsAccount.Range("A2").Value = "Joe Smith"
$wsAccount.Range("B2").Value = "Dollar"

[void]$wsAccount.Range("B3").Select()

# openAccountTest becomes the active sheet/tab
[void]$wsAccount.Activate()
